# Tripadvisor New Orleans shard 182 - update:
#   1. hotel_info: insert a new "State" column right after "Hotel_Name"
#      (before "City"), populated with "Louisiana" for the data row.
#   2. Reorder the worksheet tabs so "review_info" comes before "hotel_info".

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new State/Louisiana column into hotel_info -------------
$hotel = $wb.Worksheets.Item("hotel_info")

# Column C currently holds "City" (header) / "New Orleans" (data).
# Inserting here shifts City/Zip/TA_ReviewURL/... one column to the right
# and leaves a blank column C for the new field.
$hotel.Columns.Item(3).Insert()

$hotel.Range("C1").Value = "State"
$hotel.Range("C2").Value = "Louisiana"

# --- 2. Put review_info ahead of hotel_info in the tab order --------------
$review = $wb.Worksheets.Item("review_info")
$review.Move($hotel)
